$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns to match the latest scrape.
# Price values that are purely numeric-looking (e.g. "608.74") need to be
# force-written as text (matching the source data, which stores every price as
# a literal string, including values like "1.00" that Excel would otherwise
# coerce to the number 1). We flip the cell to Text format, assign the literal
# string, then restore the default "Normal" style so no stray formatting lingers.

$ws.Range("D2").Value = '66.780.10'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '3.256.77'
$ws.Range("E3").Value = '  +2.69%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.77%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.256.95'
$ws.Range("E8").Value = '  +2.71%  '

$ws.Range("E9").Value = '  -0.50%  '

$ws.Range("E10").Value = '  +2.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.98%  '

$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.89%  '

$ws.Range("D15").Value = '3.789.17'
$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("D16").Value = '66.781.40'
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.67%  '

$ws.Range("D18").Value = '3.265.50'
$ws.Range("E18").Value = '  +2.97%  '

$ws.Range("E19").Value = '  +1.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '509.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.21%  '

$ws.Range("E21").Value = '  +0.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.753'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.16%  '

$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.43%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  +59.20%  '

$ws.Range("E28").Value = '  +1.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.73%  '

$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("E35").Value = '  -3.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.39'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +22.25%  '

$ws.Range("E38").Value = '  +18.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '496.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0428'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.07%  '

$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.295'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.73%  '

$ws.Range("D46").Value = '2.973.42'
$ws.Range("E46").Value = '  +4.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.45%  '

$ws.Range("E48").Value = '  +5.39%  '

$ws.Range("E49").Value = '  +3.05%  '

$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("E51").Value = '  -1.27%  '

